$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("C3").Value = 155869
$ws.Range("C4").Value = 146968
$ws.Range("C5").Value = 8901
$ws.Range("C8").Value = 63.79
